$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Parent company" (B) and "Location County/City" (E) columns are no
# longer needed -- delete them entirely, shifting remaining columns left.
$ws.Range("E:E").Delete() | Out-Null
$ws.Range("B:B").Delete() | Out-Null

# Move the active selection (matches the author's final cursor position).
$ws.Range("E8").Select() | Out-Null
